# Update the cryptocurrency listing with the latest scraped values/percentages,
# and restore the row ordering that the upstream scraper produced this run
# (a handful of adjacent rows swapped rank position).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is purely numeric-looking (e.g. "530.77") must be
# force-formatted as Text first, otherwise Excel auto-converts them to numbers.
$textCells = @("D5", "D6", "D8", "D14", "D17", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D28", "D31", "D33", "D34", "D35", "D36", "D38", "D40", "D42", "D43", "D44", "D45", "D49", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "58.656.76"
$ws.Range("E2").Value = "  +1.01%  "
# Row 3 - Ethereum
$ws.Range("D3").Value = "3.150.83"
$ws.Range("E3").Value = "  +0.43%  "
# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.02%  "
# Row 5 - BNB
$ws.Range("D5").Value = "530.77"
$ws.Range("E5").Value = "  -0.41%  "
# Row 6 - Solana
$ws.Range("D6").Value = "139.69"
$ws.Range("E6").Value = "  +1.02%  "
# Row 7 - USDC
$ws.Range("E7").Value = "  +0.04%  "
# Row 8 - XRP
$ws.Range("D8").Value = "0.533"
$ws.Range("E8").Value = "  +14.20%  "
# Row 9 - Toncoin
$ws.Range("E9").Value = "  +0.17%  "
# Row 10 - Cardano
$ws.Range("E10").Value = "  +4.85%  "
# Row 11 - Dogecoin
$ws.Range("E11").Value = "  +2.95%  "
# Row 12 - TRON
$ws.Range("E12").Value = "  +2.72%  "
# Row 13 - WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "3.694.80"
$ws.Range("E13").Value = "  +0.67%  "
# Row 14 - Avalanche
$ws.Range("D14").Value = "25.97"
$ws.Range("E14").Value = "  +0.99%  "
# Row 15 - ShibaInu
$ws.Range("E15").Value = "  +4.32%  "
# Row 16 - WrappedBTC
$ws.Range("D16").Value = "58.690.09"
$ws.Range("E16").Value = "  +0.99%  "
# Row 17 - Polkadot
$ws.Range("D17").Value = "6.24"
$ws.Range("E17").Value = "  +3.68%  "
# Row 18 - WrappedEther
$ws.Range("D18").Value = "3.145.22"
$ws.Range("E18").Value = "  +0.32%  "
# Row 19 - Chainlink
$ws.Range("D19").Value = "13.02"
$ws.Range("E19").Value = "  +2.55%  "
# Row 20 - Uniswap
$ws.Range("D20").Value = "8.14"
$ws.Range("E20").Value = "  +0.09%  "
# Row 21 - BitcoinCash
$ws.Range("D21").Value = "372.02"
$ws.Range("E21").Value = "  +3.94%  "
# Row 22 - LEO
$ws.Range("D22").Value = "5.80"
$ws.Range("E22").Value = "  +1.56%  "
# Row 23 - Dai
$ws.Range("E23").Value = "  +0.27%  "
# Row 24 - Litecoin -> Polygon (swap)
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").Value = "0.523"
$ws.Range("E24").Value = "  +3.64%  "
# Row 25 - Polygon -> Litecoin (swap)
$ws.Range("B25").Value = "Litecoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D25").Value = "69.80"
$ws.Range("E25").Value = "  +0.96%  "
# Row 26 - Kaspa
$ws.Range("E26").Value = "  -0.09%  "
# Row 27 - Binance-PegBSC-USD
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").Value = "  -0.27%  "
# Row 28 - InternetComputer(DFINITY)
$ws.Range("D28").Value = "8.27"
# Row 29 - PEPE
$ws.Range("D29").Value = "0.0₃0860"
$ws.Range("E29").Value = "  -1.76%  "
# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.16%  "
# Row 31 - EthereumClassic
$ws.Range("D31").Value = "22.10"
$ws.Range("E31").Value = "  +3.04%  "
# Row 32 - RenderToken
$ws.Range("E32").Value = "  -1.19%  "
# Row 33 - NEARProtocol
$ws.Range("D33").Value = "5.17"
$ws.Range("E33").Value = "  +3.17%  "
# Row 34 - Fetch.AI
$ws.Range("D34").Value = "1.16"
$ws.Range("E34").Value = "  +1.66%  "
# Row 35 - Monero
$ws.Range("D35").Value = "158.76"
$ws.Range("E35").Value = "  -0.01%  "
# Row 36 - Aptos
$ws.Range("D36").Value = "6.27"
$ws.Range("E36").Value = "  +3.17%  "
# Row 37 - ImmutableX
$ws.Range("E37").Value = "  +5.62%  "
# Row 38 - EnergySwap
$ws.Range("D38").Value = "25.17"
$ws.Range("E38").Value = "  -2.50%  "
# Row 39 - Stacks
$ws.Range("E39").Value = "  -1.19%  "
# Row 40 - Maker -> Hedera (swap)
$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").Value = "0.0685"
$ws.Range("E40").Value = "  +2.36%  "
# Row 41 - Hedera -> Maker (swap)
$ws.Range("B41").Value = "Maker"
$ws.Range("C41").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D41").Value = "2.632.60"
$ws.Range("E41").Value = "  +4.99%  "
# Row 42 - Filecoin
$ws.Range("D42").Value = "4.25"
$ws.Range("E42").Value = "  +6.47%  "
# Row 43 - OKB
$ws.Range("D43").Value = "38.91"
$ws.Range("E43").Value = "  +3.73%  "
# Row 44 - Mantle -> VeChain (swap)
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0286"
$ws.Range("E44").Value = "  +6.17%  "
# Row 45 - VeChain -> Mantle (swap)
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "0.709"
$ws.Range("E45").Value = "  +1.03%  "
# Row 46 - FirstDigitalUSD
$ws.Range("E46").Value = "  +0.02%  "
# Row 47 - RenzoRestakedETH
$ws.Range("D47").Value = "3.195.25"
$ws.Range("E47").Value = "  +0.59%  "
# Row 48 - Stellar
$ws.Range("E48").Value = "  +13.50%  "
# Row 49 - ONDO -> Cosmos (swap)
$ws.Range("B49").Value = "Cosmos"
$ws.Range("C49").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D49").Value = "6.20"
$ws.Range("E49").Value = "  +2.62%  "
# Row 50 - Cosmos -> ONDO (swap)
$ws.Range("B50").Value = "ONDO"
$ws.Range("C50").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D50").Value = "0.979"
$ws.Range("E50").Value = "  +0.08%  "
# Row 51 - InjectiveProtocol
$ws.Range("D51").Value = "20.33"
$ws.Range("E51").Value = "  +2.69%  "
